$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.8986206666666666
$ws.Range("H2").Value = 2.695862
$ws.Range("I2").Value = 0.1661804693926261
$ws.Range("J2").Value = 0.1661804693926262
$ws.Range("M2").Value = 4.230734666666667
$ws.Range("N2").Value = 12.692204
$ws.Range("O2").Value = 0.3081346507358854
$ws.Range("P2").Value = 0.3081346507358855
$ws.Range("Q2").Value = 3.801825606649778
$ws.Range("R2").Value = 34.216430459848
$ws.Range("S2").Value = 0.05120596089542235
$ws.Range("T2").Value = 0.05120596089542236

$ws.Range("G3").Value = 0.8986206666666666
$ws.Range("H3").Value = 2.695862
$ws.Range("I3").Value = 0.1661804693926261
$ws.Range("J3").Value = 0.1661804693926262
$ws.Range("O3").Value = 0.6225996527787135
$ws.Range("P3").Value = 0.6225996527787135
$ws.Range("Q3").Value = 7.681756326244
$ws.Range("R3").Value = 69.135806936196
$ws.Range("S3").Value = 0.1034639025424527
$ws.Range("T3").Value = 0.1034639025424527

$ws.Range("G4").Value = 0.8986206666666666
$ws.Range("H4").Value = 2.695862
$ws.Range("I4").Value = 0.1661804693926261
$ws.Range("J4").Value = 0.1661804693926262
$ws.Range("M4").Value = 0.9510283333333334
$ws.Range("N4").Value = 2.853085
$ws.Range("O4").Value = 0.0692656964854011
$ws.Range("P4").Value = 0.0692656964854011
$ws.Range("Q4").Value = 0.8546137149188888
$ws.Range("R4").Value = 7.691523434270001
$ws.Range("S4").Value = 0.01151060595475113
$ws.Range("T4").Value = 0.01151060595475113

$ws.Range("I5").Value = 0.3901029163453022
$ws.Range("J5").Value = 0.3901029163453023
$ws.Range("M5").Value = 4.230734666666667
$ws.Range("N5").Value = 12.692204
$ws.Range("O5").Value = 0.3081346507358854
$ws.Range("P5").Value = 0.3081346507358855
$ws.Range("Q5").Value = 8.924654395374667
$ws.Range("R5").Value = 80.321889558372
$ws.Range("S5").Value = 0.12020422587911
$ws.Range("T5").Value = 0.1202042258791101

$ws.Range("I6").Value = 0.3901029163453022
$ws.Range("J6").Value = 0.3901029163453023
$ws.Range("O6").Value = 0.6225996527787135
$ws.Range("P6").Value = 0.6225996527787135
$ws.Range("S6").Value = 0.2428779402645487
$ws.Range("T6").Value = 0.2428779402645487

$ws.Range("I7").Value = 0.3901029163453022
$ws.Range("J7").Value = 0.3901029163453023
$ws.Range("M7").Value = 0.9510283333333334
$ws.Range("N7").Value = 2.853085
$ws.Range("O7").Value = 0.0692656964854011
$ws.Range("P7").Value = 0.0692656964854011
$ws.Range("Q7").Value = 2.006176199628333
$ws.Range("R7").Value = 18.055585796655
$ws.Range("S7").Value = 0.02702075020164352
$ws.Range("T7").Value = 0.02702075020164352

$ws.Range("G8").Value = 2.399397
$ws.Range("H8").Value = 7.198191
$ws.Range("I8").Value = 0.4437166142620716
$ws.Range("J8").Value = 0.4437166142620716
$ws.Range("M8").Value = 4.230734666666667
$ws.Range("N8").Value = 12.692204
$ws.Range("O8").Value = 0.3081346507358854
$ws.Range("P8").Value = 0.3081346507358855
$ws.Range("Q8").Value = 10.151212066996
$ws.Range("R8").Value = 91.360908602964
$ws.Range("S8").Value = 0.136724463961353
$ws.Range("T8").Value = 0.136724463961353

$ws.Range("G9").Value = 2.399397
$ws.Range("H9").Value = 7.198191
$ws.Range("I9").Value = 0.4437166142620716
$ws.Range("J9").Value = 0.4437166142620716
$ws.Range("O9").Value = 0.6225996527787135
$ws.Range("P9").Value = 0.6225996527787135
$ws.Range("Q9").Value = 20.510971723242
$ws.Range("R9").Value = 184.598745509178
$ws.Range("S9").Value = 0.2762578099717122
$ws.Range("T9").Value = 0.2762578099717122

$ws.Range("G10").Value = 2.399397
$ws.Range("H10").Value = 7.198191
$ws.Range("I10").Value = 0.4437166142620716
$ws.Range("J10").Value = 0.4437166142620716
$ws.Range("M10").Value = 0.9510283333333334
$ws.Range("N10").Value = 2.853085
$ws.Range("O10").Value = 0.0692656964854011
$ws.Range("P10").Value = 0.0692656964854011
$ws.Range("Q10").Value = 2.281894529915
$ws.Range("R10").Value = 20.537050769235
$ws.Range("S10").Value = 0.03073434032900645
$ws.Range("T10").Value = 0.03073434032900645

